# Atualiza granularidade da consulta avançada
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# The "Anexo" (column I) granularity was missing an "x" mark for rows 2-10;
# bring it in line with the neighboring H/J columns of the advanced query grid.
for ($r = 2; $r -le 10; $r++) {
    $ws1.Cells.Item($r, 9).Value = "x"
}

# Update Sheet3's lingering selection first (selecting on it would otherwise
# make it the active tab again if done after Sheet1 is activated).
$ws3.Range("A92").Select()

# Sheet1 becomes the active/visible tab, with the cursor resting on I11.
$ws1.Activate()
$ws1.Range("I11").Select()

$wb.Save()
